# Auto-generated script to apply numeric corrections to Sheets per diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 633.4706
$ws.Cells.Item(98, 9).Value = 650.86664
$ws.Cells.Item(98, 11).Value = 650.86664
$ws.Cells.Item(98, 13).Value = 847.13336
$ws.Cells.Item(116, 8).Value = 1890.1111
$ws.Cells.Item(116, 9).Value = 1751
$ws.Cells.Item(116, 11).Value = 1751
$ws.Cells.Item(116, 13).Value = 1691
$ws.Cells.Item(122, 8).Value = 633.4706
$ws.Cells.Item(122, 9).Value = 650.86664
$ws.Cells.Item(122, 11).Value = 1952.59992
$ws.Cells.Item(122, 13).Value = 497.4000800000001
$ws.Cells.Item(132, 8).Value = 1486562
$ws.Cells.Item(132, 9).Value = 1751.8334
$ws.Cells.Item(132, 10).Value = 16334663
$ws.Cells.Item(132, 11).Value = 5255.5002
$ws.Cells.Item(132, 12).Value = 49003989
$ws.Cells.Item(132, 13).Value = -2725.5002
$ws.Cells.Item(132, 14).Value = -49009049
$ws.Cells.Item(138, 8).Value = 2198.6858
$ws.Cells.Item(138, 9).Value = 973.3
$ws.Cells.Item(138, 10).Value = 2688.84
$ws.Cells.Item(138, 11).Value = 2919.9
$ws.Cells.Item(138, 12).Value = 8066.52
$ws.Cells.Item(138, 13).Value = 2220.1
$ws.Cells.Item(138, 14).Value = -18346.52

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1828.5
$ws.Cells.Item(2, 9).Value = 1442.8667
$ws.Cells.Item(2, 10).Value = 2273.4614
$ws.Cells.Item(2, 11).Value = 1442.8667
$ws.Cells.Item(2, 12).Value = 2273.4614
$ws.Cells.Item(2, 13).Value = -1329.8667
$ws.Cells.Item(2, 14).Value = -2499.4614
$ws.Cells.Item(45, 8).Value = 2909
$ws.Cells.Item(45, 9).Value = 3054.7693
$ws.Cells.Item(45, 10).Value = 1014
$ws.Cells.Item(45, 11).Value = 3054.7693
$ws.Cells.Item(45, 12).Value = 1014
$ws.Cells.Item(45, 13).Value = -2677.7693
$ws.Cells.Item(45, 14).Value = -1768
$ws.Cells.Item(74, 8).Value = 8692369
$ws.Cells.Item(74, 9).Value = 12552400
$ws.Cells.Item(74, 10).Value = 114521.11
$ws.Cells.Item(74, 11).Value = 12552400
$ws.Cells.Item(74, 12).Value = 114521.11
$ws.Cells.Item(74, 13).Value = -12551526
$ws.Cells.Item(74, 14).Value = -116269.11
$ws.Cells.Item(77, 8).Value = 8692369
$ws.Cells.Item(77, 9).Value = 12552400
$ws.Cells.Item(77, 10).Value = 114521.11
$ws.Cells.Item(77, 11).Value = 62762000
$ws.Cells.Item(77, 12).Value = 572605.55
$ws.Cells.Item(77, 13).Value = -62757632
$ws.Cells.Item(77, 14).Value = -581341.55
$ws.Cells.Item(110, 8).Value = 1243
$ws.Cells.Item(110, 9).Value = 650.1818
$ws.Cells.Item(110, 10).Value = 2058.125
$ws.Cells.Item(110, 11).Value = 650.1818
$ws.Cells.Item(110, 12).Value = 2058.125
$ws.Cells.Item(110, 13).Value = 1394.8182
$ws.Cells.Item(110, 14).Value = -6148.125
$ws.Cells.Item(116, 8).Value = 1828.5
$ws.Cells.Item(116, 9).Value = 1442.8667
$ws.Cells.Item(116, 10).Value = 2273.4614
$ws.Cells.Item(116, 11).Value = 1442.8667
$ws.Cells.Item(116, 12).Value = 2273.4614
$ws.Cells.Item(116, 13).Value = 851.1333
$ws.Cells.Item(116, 14).Value = -6861.4614
$ws.Cells.Item(122, 8).Value = 22224220
$ws.Cells.Item(122, 9).Value = 2663.3333
$ws.Cells.Item(122, 10).Value = 55556556
$ws.Cells.Item(122, 11).Value = 7989.999899999999
$ws.Cells.Item(122, 12).Value = 166669668
$ws.Cells.Item(122, 13).Value = -5539.999899999999
$ws.Cells.Item(122, 14).Value = -166674568
$ws.Cells.Item(139, 8).Value = 40045.312
$ws.Cells.Item(139, 10).Value = 40045.312
$ws.Cells.Item(139, 12).Value = 40045.312
$ws.Cells.Item(139, 14).Value = -50325.312

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1828.5
$ws.Cells.Item(3, 9).Value = 1442.8667
$ws.Cells.Item(3, 10).Value = 2273.4614
$ws.Cells.Item(3, 11).Value = 1442.8667
$ws.Cells.Item(3, 12).Value = 2273.4614
$ws.Cells.Item(3, 13).Value = -1328.8667
$ws.Cells.Item(3, 14).Value = -2501.4614
$ws.Cells.Item(20, 8).Value = 749.5
$ws.Cells.Item(20, 9).Value = 552.5454999999999
$ws.Cells.Item(20, 10).Value = 1110.5834
$ws.Cells.Item(20, 11).Value = 552.5454999999999
$ws.Cells.Item(20, 12).Value = 1110.5834
$ws.Cells.Item(20, 13).Value = -305.5454999999999
$ws.Cells.Item(20, 14).Value = -1604.5834
$ws.Cells.Item(99, 8).Value = 1029.3529
$ws.Cells.Item(99, 9).Value = 1060
$ws.Cells.Item(99, 10).Value = 985.5714
$ws.Cells.Item(99, 11).Value = 1060
$ws.Cells.Item(99, 12).Value = 985.5714
$ws.Cells.Item(99, 13).Value = 438
$ws.Cells.Item(99, 14).Value = -3981.5714
$ws.Cells.Item(134, 8).Value = 3391
$ws.Cells.Item(134, 9).Value = 2671.158
$ws.Cells.Item(134, 10).Value = 5344.857
$ws.Cells.Item(134, 11).Value = 8013.474
$ws.Cells.Item(134, 12).Value = 16034.571
$ws.Cells.Item(134, 13).Value = -5478.474
$ws.Cells.Item(134, 14).Value = -21104.571

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1753.2667
$ws.Cells.Item(31, 9).Value = 1204.4762
$ws.Cells.Item(31, 10).Value = 3033.7778
$ws.Cells.Item(31, 11).Value = 1204.4762
$ws.Cells.Item(31, 12).Value = 3033.7778
$ws.Cells.Item(31, 13).Value = -909.4762000000001
$ws.Cells.Item(31, 14).Value = -3623.7778
$ws.Cells.Item(34, 8).Value = 1753.2667
$ws.Cells.Item(34, 9).Value = 1204.4762
$ws.Cells.Item(34, 10).Value = 3033.7778
$ws.Cells.Item(34, 11).Value = 1204.4762
$ws.Cells.Item(34, 12).Value = 3033.7778
$ws.Cells.Item(34, 13).Value = -1002.4762
$ws.Cells.Item(34, 14).Value = -3437.7778
$ws.Cells.Item(107, 8).Value = 466.42856
$ws.Cells.Item(107, 9).Value = 426.66666
$ws.Cells.Item(107, 10).Value = 538
$ws.Cells.Item(107, 11).Value = 426.66666
$ws.Cells.Item(107, 12).Value = 538
$ws.Cells.Item(107, 13).Value = 1493.33334
$ws.Cells.Item(107, 14).Value = -4378
$ws.Cells.Item(132, 8).Value = 34598.844
$ws.Cells.Item(132, 9).Value = 3119.8333
$ws.Cells.Item(132, 10).Value = 75071.86
$ws.Cells.Item(132, 11).Value = 9359.499899999999
$ws.Cells.Item(132, 12).Value = 225215.58
$ws.Cells.Item(132, 13).Value = -6829.499899999999
$ws.Cells.Item(132, 14).Value = -230275.58
$ws.Cells.Item(134, 8).Value = 33669
$ws.Cells.Item(134, 9).Value = 2278.3794
$ws.Cells.Item(134, 10).Value = 163715.86
$ws.Cells.Item(134, 11).Value = 6835.138199999999
$ws.Cells.Item(134, 12).Value = 491147.58
$ws.Cells.Item(134, 13).Value = -4300.138199999999
$ws.Cells.Item(134, 14).Value = -496217.58

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 700
$ws.Cells.Item(2, 9).Value = 700
$ws.Cells.Item(2, 11).Value = 4200
$ws.Cells.Item(2, 13).Value = -4087
$ws.Cells.Item(22, 8).Value = 2862.25
$ws.Cells.Item(22, 9).Value = 449
$ws.Cells.Item(22, 10).Value = 3666.6667
$ws.Cells.Item(22, 11).Value = 1347
$ws.Cells.Item(22, 12).Value = 11000.0001
$ws.Cells.Item(22, 13).Value = -1178
$ws.Cells.Item(22, 14).Value = -11338.0001
$ws.Cells.Item(27, 8).Value = 2862.25
$ws.Cells.Item(27, 9).Value = 449
$ws.Cells.Item(27, 10).Value = 3666.6667
$ws.Cells.Item(27, 11).Value = 1347
$ws.Cells.Item(27, 12).Value = 11000.0001
$ws.Cells.Item(27, 13).Value = -1245
$ws.Cells.Item(27, 14).Value = -11204.0001
$ws.Cells.Item(76, 8).Value = 3347.3076
$ws.Cells.Item(76, 9).Value = 1000
$ws.Cells.Item(76, 10).Value = 3542.9167
$ws.Cells.Item(76, 11).Value = 3000
$ws.Cells.Item(76, 12).Value = 10628.7501
$ws.Cells.Item(76, 13).Value = -2617
$ws.Cells.Item(76, 14).Value = -11394.7501
$ws.Cells.Item(79, 8).Value = 3347.3076
$ws.Cells.Item(79, 9).Value = 1000
$ws.Cells.Item(79, 10).Value = 3542.9167
$ws.Cells.Item(79, 11).Value = 3000
$ws.Cells.Item(79, 12).Value = 10628.7501
$ws.Cells.Item(79, 13).Value = -1674
$ws.Cells.Item(79, 14).Value = -13280.7501
$ws.Cells.Item(131, 8).Value = 943.2273
$ws.Cells.Item(131, 9).Value = 534.9
$ws.Cells.Item(131, 10).Value = 1016.1429
$ws.Cells.Item(131, 11).Value = 1604.7
$ws.Cells.Item(131, 12).Value = 3048.4287
$ws.Cells.Item(131, 13).Value = 3435.3
$ws.Cells.Item(131, 14).Value = -13128.4287
$ws.Cells.Item(132, 8).Value = 835.4
$ws.Cells.Item(132, 9).Value = 700
$ws.Cells.Item(132, 10).Value = 869.25
$ws.Cells.Item(132, 11).Value = 6300
$ws.Cells.Item(132, 12).Value = 7823.25
$ws.Cells.Item(132, 13).Value = -3770
$ws.Cells.Item(132, 14).Value = -12883.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(118, 8).Value = 50310
$ws.Cells.Item(118, 10).Value = 50310
$ws.Cells.Item(118, 12).Value = 50310
$ws.Cells.Item(118, 14).Value = -53624
$ws.Cells.Item(122, 8).Value = 2660.3
$ws.Cells.Item(122, 9).Value = 1765
$ws.Cells.Item(122, 10).Value = 3044
$ws.Cells.Item(122, 11).Value = 5295
$ws.Cells.Item(122, 12).Value = 9132
$ws.Cells.Item(122, 13).Value = -2845
$ws.Cells.Item(122, 14).Value = -14032
$ws.Cells.Item(139, 8).Value = 81318
$ws.Cells.Item(139, 10).Value = 81318
$ws.Cells.Item(139, 12).Value = 81318
$ws.Cells.Item(139, 14).Value = -91598

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 77929.5
$ws.Cells.Item(132, 9).Value = 38442.285
$ws.Cells.Item(132, 10).Value = 170066.33
$ws.Cells.Item(132, 11).Value = 115326.855
$ws.Cells.Item(132, 12).Value = 510198.99
$ws.Cells.Item(132, 13).Value = -112796.855
$ws.Cells.Item(132, 14).Value = -515258.99
$ws.Cells.Item(136, 8).Value = 251066.88
$ws.Cells.Item(136, 9).Value = 200907
$ws.Cells.Item(136, 10).Value = 334666.66
$ws.Cells.Item(136, 11).Value = 602721
$ws.Cells.Item(136, 12).Value = 1003999.98
$ws.Cells.Item(136, 13).Value = -600171
$ws.Cells.Item(136, 14).Value = -1009099.98

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 2390.8667
$ws.Cells.Item(122, 9).Value = 972.875
$ws.Cells.Item(122, 11).Value = 2918.625
$ws.Cells.Item(122, 13).Value = -468.625
$ws.Cells.Item(136, 8).Value = 39440.117
$ws.Cells.Item(136, 9).Value = 21324.203
$ws.Cells.Item(136, 11).Value = 63972.609
$ws.Cells.Item(136, 13).Value = -61422.609
